$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "VISCERALGINE 10MG/5ML SYRUP 120 ML" line item (row 113) entirely.
# Deleting the row shifts everything below it up by one.
$ws.Rows("113:113").Delete()

# The grand-total cell (now on row 142 after the shift) is a hard-coded value,
# not a formula, so it has to be corrected by hand: subtract the removed
# item's amount (35.00) from the old total.
$ws.Range("P142").Value = 7646.79
